$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; this shifts existing rows 26-40 down to 27-41
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with its data.
# (Same as the values that were in the old row 26 before the shift, except for
# the Fecha, Variedad, Volumen and Origen fields which carry the new data.)
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value = "Ñuble"
$ws.Range("D26").Value = 44529
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = 100112031
$ws.Range("G26").Value = "Poroto verde"
$ws.Range("H26").Value = "Magnum"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 60
$ws.Range("K26").Value = 29000
$ws.Range("L26").Value = 30000
$ws.Range("M26").Value = 29500
$ws.Range("N26").Value = "$/saco 25 kilos"
$ws.Range("O26").Value = "Región Metropolitana"
$ws.Range("P26").Value = 1180
$ws.Range("Q26").Value = 25
$ws.Range("R26").Value = "Hortaliza"
